$p = $ppt.ActivePresentation

# The deck currently ends: ... , "References" (slide 7).
# Insert a new "Title and Content" slide at the end, then move it so it
# lands right before "References" -- i.e. it becomes the new slide 7 and
# "References" is pushed down to slide 8.
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)
[void]$newSlide.MoveTo(7)

$newSlide = $p.Slides.Item(7)

# Title placeholder
$newSlide.Shapes.Placeholders.Item(1).TextFrame.TextRange.Text = "Final Thoughts"

# Body / content placeholder - three bullet paragraphs
$body = $newSlide.Shapes.Placeholders.Item(2).TextFrame.TextRange
$body.Text = "These models were effective in determining the outcomes of the 4 awards I chose."
[void]$body.InsertAfter("`rTo improve it, I would try to look at a player’s progression over multiple seasons to forecast the future performance.")
[void]$body.InsertAfter("`rIt could also be useful to add in data on team success, as that would have some effects on the model. ")
